$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Create the empty separator row (row 12), mirroring the spacing used for
#    row 3 (the thick-bordered separator above the first table).
# ---------------------------------------------------------------------------
$ws.Rows.Item(12).RowHeight = 15

# ---------------------------------------------------------------------------
# 2) Duplicate the formatting of the first table (rows 4-11) down onto the
#    new second table (rows 13-21). Row 8's format is reused twice (17 & 18)
#    because the new table has one extra "Activity" line compared to the
#    first table.
# ---------------------------------------------------------------------------
$ws.Range("A4:E8").Copy()
$ws.Range("A13:E17").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("B8:C8").Copy()
$ws.Range("B18:C18").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("A9:E11").Copy()
$ws.Range("A19:E21").PasteSpecial(-4122)   # xlPasteFormats

$ws.Application.CutCopyMode = 0

# Remove the stray blank cells that PasteSpecial created for columns that
# never had any content/formatting in the source rows (6, 7 and 8).
$ws.Range("A15").Clear()
$ws.Range("D15").Clear()
$ws.Range("A16").Clear()
$ws.Range("D16").Clear()
$ws.Range("A17").Clear()
$ws.Range("D17").Clear()
$ws.Range("E17").Clear()

# ---------------------------------------------------------------------------
# 3) Fill in the values for the new "Autumn 1400" (پاییز 1400) table.
# ---------------------------------------------------------------------------
$ws.Range("A14").Value = "پاییز 1400"
$ws.Range("B14").Value = "Activity"
$ws.Range("C14").Value = "Hours"
$ws.Range("E14").Value = "Tasks Done"

$ws.Range("B15").Value = "* System Preparation"
$ws.Range("C15").Value = 6
$ws.Range("E15").Value = "• Lung/Airway Segmentation Data"

$ws.Range("B16").Value = "* Data Preparation"
$ws.Range("C16").Value = 5
$ws.Range("E16").Value = "• Deep Segmentation Model"

$ws.Range("B17").Value = "* Deep Model Code"
$ws.Range("C17").Value = 6

$ws.Range("B18").Value = "* Meetings"
$ws.Range("C18").Value = 3

$ws.Range("B19").Value = "• Total Hours"
$ws.Range("C19").Formula = "=SUM(C15:C18)"

$ws.Range("C20").Value = "@Parsiss"
$ws.Range("D20").Value = 7

$ws.Range("C21").Value = "@Home"
$ws.Range("D21").Formula = "=C19-D20"

# ---------------------------------------------------------------------------
# 4) Match the saved selection/active cell of the workbook.
# ---------------------------------------------------------------------------
$ws.Range("E24").Select()
